# Update the "Canceled" value for the transaction in row 3 (Writing essays /
# Essay for friend 6) from "No" to "Yes" on the Statement sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Statement")

$ws.Range("F3").Value = "Yes"
